# Add a new worksheet "header" after "icecream", positioned as the active (last) tab.
# It re-presents the icecream table with the header row split across two rows
# (A/C/E/G on row 1, B/D/F on row 2), matching the new "header" layout.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# icecream sheet: selection moves from the old G12 to the sorted data block A1:G10.
$ws1.Range("A1:G10").Select()

$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "header"

$new.Range("A1").Value2 = "Type"
$new.Range("A1").Font.Bold = $true
$new.Range("C1").Value2 = "June"
$new.Range("C1").Font.Bold = $true
$new.Range("E1").Value2 = "August"
$new.Range("E1").Font.Bold = $true
$new.Range("G1").Value2 = "Profit"
$new.Range("G1").Font.Bold = $true
$new.Range("B2").Value2 = "Country"
$new.Range("B2").Font.Bold = $true
$new.Range("D2").Value2 = "July"
$new.Range("D2").Font.Bold = $true
$new.Range("F2").Value2 = "Total"
$new.Range("F2").Font.Bold = $true
$new.Range("A3").Value2 = "Banana"
$new.Range("A3").Font.Bold = $true
$new.Range("B3").Value2 = "BE"
$new.Range("C3").Value2 = 170
$new.Range("D3").Value2 = 690
$new.Range("E3").Value2 = 520
$new.Range("F3").Value2 = 1380
$new.Range("G3").Value2 = "YES"
$new.Range("B4").Value2 = "DE"
$new.Range("C4").Value2 = 610
$new.Range("D4").Value2 = 640
$new.Range("E4").Value2 = 320
$new.Range("F4").Value2 = 1570
$new.Range("G4").Value2 = "NO"
$new.Range("B5").Value2 = "DE"
$new.Range("C5").Value2 = 250
$new.Range("D5").Value2 = 650
$new.Range("E5").Value2 = 630
$new.Range("F5").Value2 = 1530
$new.Range("G5").Value2 = "YES"
$new.Range("A6").Value2 = "Chocolate"
$new.Range("A6").Font.Bold = $true
$new.Range("B6").Value2 = "BE"
$new.Range("C6").Value2 = 560
$new.Range("D6").Value2 = 320
$new.Range("E6").Value2 = 140
$new.Range("F6").Value2 = 1020
$new.Range("G6").Value2 = "YES"
$new.Range("B7").Value2 = "FR"
$new.Range("C7").Value2 = 430
$new.Range("D7").Value2 = 350
$new.Range("E7").Value2 = 300
$new.Range("F7").Value2 = 1080
$new.Range("G7").Value2 = "YES"
$new.Range("C8").Value2 = "N/A"
$new.Range("D8").Value2 = "N/A"
$new.Range("E8").Value2 = "N/A"
$new.Range("F8").Value2 = "N/A"
$new.Range("B9").Value2 = "NL"
$new.Range("C9").Value2 = 210
$new.Range("D9").Value2 = 280
$new.Range("E9").Value2 = 270
$new.Range("F9").Value2 = 760
$new.Range("G9").Value2 = "NO"
$new.Range("A10").Value2 = "Speculaas"
$new.Range("A10").Font.Bold = $true
$new.Range("B10").Value2 = "BE"
$new.Range("C10").Value2 = 300
$new.Range("D10").Value2 = 270
$new.Range("E10").Value2 = 290
$new.Range("F10").Value2 = 860
$new.Range("G10").Value2 = "NO"
$new.Range("A11").Value2 = "Vanilla"
$new.Range("A11").Font.Bold = $true
$new.Range("B11").Value2 = "BE"
$new.Range("C11").Value2 = 610
$new.Range("D11").Value2 = 190
$new.Range("E11").Value2 = 670
$new.Range("F11").Value2 = 1470
$new.Range("G11").Value2 = "YES"

$new.Range("F2").Select()
